# Summer week 11 inputs - append new matchup rows to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(4,2,6,1),
    @(5,0,4,3),
    @(4,0,6,2),
    @(3,3,2,0),
    @(3,2,3,1),
    @(6,2,5,1),
    @(4,2,3,1),
    @(6,0,5,2),
    @(4,1,4,2),
    @(3,2,3,1),
    @(7,2,7,0),
    @(4,2,5,0),
    @(5,2,5,0),
    @(3,3,2,0),
    @(4,2,2,1),
    @(3,0,3,3),
    @(5,0,5,2),
    @(5,1,6,2),
    @(6,3,4,0),
    @(4,2,4,1)
)

$startRow = 1110
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update the visible window / selection to match the scrolled state after entry
$ws.Activate()
$excel.ActiveWindow.ScrollRow = $startRow
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J1124").Select()
